$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2491913.2
$ws.Range("J17").Value = 2587723.2
$ws.Range("L17").Value = 7763169.600000001
$ws.Range("N17").Value = -7763505.600000001
$ws.Range("H33").Value = 22561.8
$ws.Range("J33").Value = 4667.3335
$ws.Range("L33").Value = 4667.3335
$ws.Range("N33").Value = -5125.3335
$ws.Range("H51").Value = 4976.952
$ws.Range("I51").Value = 4841.3335
$ws.Range("J51").Value = 4999.5557
$ws.Range("K51").Value = 4841.3335
$ws.Range("L51").Value = 4999.5557
$ws.Range("M51").Value = -4357.3335
$ws.Range("N51").Value = -5967.5557
$ws.Range("H53").Value = 691.4761999999999
$ws.Range("I53").Value = 604.8461
$ws.Range("J53").Value = 832.25
$ws.Range("K53").Value = 604.8461
$ws.Range("L53").Value = 832.25
$ws.Range("M53").Value = 32.15390000000002
$ws.Range("N53").Value = -2106.25
$ws.Range("H94").Value = 914.8333
$ws.Range("I94").Value = 914.8333
$ws.Range("K94").Value = 914.8333
$ws.Range("M94").Value = -463.8333
$ws.Range("H103").Value = 760.8889
$ws.Range("I103").Value = 750.4286
$ws.Range("J103").Value = 797.5
$ws.Range("K103").Value = 2251.2858
$ws.Range("L103").Value = 2392.5
$ws.Range("M103").Value = -1665.2858
$ws.Range("N103").Value = -3564.5
$ws.Range("H137").Value = 13427.823
$ws.Range("I137").Value = 5969.857
$ws.Range("K137").Value = 17909.571
$ws.Range("M137").Value = -15359.571

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3595.5
$ws.Range("I45").Value = 3458.1667
$ws.Range("K45").Value = 3458.1667
$ws.Range("M45").Value = -3081.1667
$ws.Range("H61").Value = 30306930
$ws.Range("I61").Value = 37040230
$ws.Range("K61").Value = 37040230
$ws.Range("M61").Value = -37040018
$ws.Range("H110").Value = 22108.47
$ws.Range("I110").Value = 26096.072
$ws.Range("K110").Value = 26096.072
$ws.Range("M110").Value = -24051.072
$ws.Range("H132").Value = 26386496
$ws.Range("I132").Value = 8260.235000000001
$ws.Range("K132").Value = 24780.705
$ws.Range("M132").Value = -22250.705
$ws.Range("H136").Value = 30306930
$ws.Range("I136").Value = 37040230
$ws.Range("K136").Value = 111120690
$ws.Range("M136").Value = -111118140

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3425.5386
$ws.Range("I99").Value = 1401.25
$ws.Range("J99").Value = 6664.4
$ws.Range("K99").Value = 1401.25
$ws.Range("L99").Value = 6664.4
$ws.Range("M99").Value = 96.75
$ws.Range("N99").Value = -9660.4
$ws.Range("H124").Value = 149999
$ws.Range("J124").Value = 149999
$ws.Range("L124").Value = 149999
$ws.Range("N124").Value = -159819

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 71434350
$ws.Range("I31").Value = 2881.1667
$ws.Range("K31").Value = 2881.1667
$ws.Range("M31").Value = -2586.1667
$ws.Range("H34").Value = 71434350
$ws.Range("I34").Value = 2881.1667
$ws.Range("K34").Value = 2881.1667
$ws.Range("M34").Value = -2679.1667
$ws.Range("H58").Value = 2602.775
$ws.Range("I58").Value = 2394.5518
$ws.Range("J58").Value = 3151.7273
$ws.Range("K58").Value = 2394.5518
$ws.Range("L58").Value = 3151.7273
$ws.Range("M58").Value = -2191.5518
$ws.Range("N58").Value = -3557.7273
$ws.Range("H99").Value = 5389.069
$ws.Range("I99").Value = 5562.6665
$ws.Range("K99").Value = 5562.6665
$ws.Range("M99").Value = -4064.6665
$ws.Range("H126").Value = 5389.069
$ws.Range("I126").Value = 5562.6665
$ws.Range("K126").Value = 16687.9995
$ws.Range("M126").Value = -14217.9995
$ws.Range("H136").Value = 2602.775
$ws.Range("I136").Value = 2394.5518
$ws.Range("J136").Value = 3151.7273
$ws.Range("K136").Value = 7183.655400000001
$ws.Range("L136").Value = 9455.1819
$ws.Range("M136").Value = -4633.655400000001
$ws.Range("N136").Value = -14555.1819

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 15152449
$ws.Range("I55").Value = 20000250
$ws.Range("J55").Value = 11112615
$ws.Range("K55").Value = 60000750
$ws.Range("L55").Value = 33337845
$ws.Range("M55").Value = -60000573
$ws.Range("N55").Value = -33338199
$ws.Range("H97").Value = 234.8
$ws.Range("I97").Value = 230
$ws.Range("J97").Value = 236
$ws.Range("K97").Value = 690
$ws.Range("L97").Value = 708
$ws.Range("N97").Value = -1700
$ws.Range("M97").Value = -194
$ws.Range("H137").Value = 2022.0714
$ws.Range("J137").Value = 1836.8334
$ws.Range("L137").Value = 5510.5002
$ws.Range("N137").Value = -15710.5002
$ws.Range("H140").Value = 3028.0715
$ws.Range("J140").Value = 3580.8572
$ws.Range("L140").Value = 10742.5716
$ws.Range("N140").Value = -21102.5716

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 23569.143
$ws.Range("I21").Value = 19989.334
$ws.Range("J21").Value = 24545.455
$ws.Range("K21").Value = 19989.334
$ws.Range("L21").Value = 24545.455
$ws.Range("M21").Value = -19816.334
$ws.Range("N21").Value = -24891.455
$ws.Range("H30").Value = 23569.143
$ws.Range("I30").Value = 19989.334
$ws.Range("J30").Value = 24545.455
$ws.Range("K30").Value = 19989.334
$ws.Range("L30").Value = 24545.455
$ws.Range("M30").Value = -19884.334
$ws.Range("N30").Value = -24755.455
$ws.Range("H70").Value = 289282.72
$ws.Range("I70").Value = 336831.84
$ws.Range("K70").Value = 336831.84
$ws.Range("M70").Value = -336561.84
$ws.Range("H73").Value = 289282.72
$ws.Range("I73").Value = 336831.84
$ws.Range("K73").Value = 336831.84
$ws.Range("M73").Value = -335895.84
$ws.Range("H80").Value = 4445.1816
$ws.Range("I80").Value = 3666.6667
$ws.Range("K80").Value = 3666.6667
$ws.Range("M80").Value = -2668.6667
$ws.Range("H83").Value = 4445.1816
$ws.Range("I83").Value = 3666.6667
$ws.Range("K83").Value = 18333.3335
$ws.Range("M83").Value = -13341.3335
$ws.Range("H102").Value = 3503.8333
$ws.Range("I102").Value = 2755.75
$ws.Range("K102").Value = 2755.75
$ws.Range("M102").Value = -1133.75
$ws.Range("H141").Value = 113786
$ws.Range("I141").Value = 84000
$ws.Range("J141").Value = 121232.5
$ws.Range("K141").Value = 84000
$ws.Range("L141").Value = 121232.5
$ws.Range("N141").Value = -131592.5
$ws.Range("M141").Value = -78820

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2501.9048
$ws.Range("I82").Value = 1104.2222
$ws.Range("K82").Value = 1104.2222
$ws.Range("M82").Value = -743.2221999999999
$ws.Range("H85").Value = 2501.9048
$ws.Range("I85").Value = 1104.2222
$ws.Range("K85").Value = 1104.2222
$ws.Range("M85").Value = 143.7778000000001
$ws.Range("H136").Value = 2892.7646
$ws.Range("I136").Value = 2246.6667
$ws.Range("K136").Value = 6740.000100000001
$ws.Range("M136").Value = -4190.000100000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 795.6667
$ws.Range("J113").Value = 832.36365
$ws.Range("L113").Value = 2497.09095
$ws.Range("N113").Value = -6837.09095
$ws.Range("H122").Value = 111224320
$ws.Range("I122").Value = 143001630
$ws.Range("K122").Value = 429004890
$ws.Range("M122").Value = -429002440
$ws.Range("H126").Value = 6123.0557
$ws.Range("I126").Value = 6548.2856
$ws.Range("K126").Value = 19644.8568
$ws.Range("M126").Value = -17174.8568
$ws.Range("H140").Value = 73926.875
$ws.Range("J140").Value = 76127
$ws.Range("L140").Value = 76127
$ws.Range("N140").Value = -86487
